$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.785.03'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.763.95'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.21'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.63'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.760.87'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +4.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.35'
$ws.Range("E11").Value = '  +2.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.49'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.394.25'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.752.27'
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.806.87'
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.32'
$ws.Range("E18").Value = '  +1.68%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.18'
$ws.Range("E20").Value = '  -1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.82'
$ws.Range("E21").Value = '  +15.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '497.97'
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.86'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000145'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.34'
$ws.Range("E26").Value = '  -1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.39'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.25'
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.55'
$ws.Range("E30").Value = '  +5.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.99'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.98'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.24'
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.912.32'
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.698.45'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.86'
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '445.01'
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.81'
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.98'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.88'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.843.52'
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.94'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("E51").Value = '  +2.17%  '
